# Applies the edit described in the commit:
#   "elaborated sanity checks. increased the proportion of new images in memory task."
# Rows 2-29 (pre-existing trials) are updated in place; rows 30-41 are newly appended
# trials. Dimension (A1:S41) is recomputed automatically by Excel from the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 362
$ws.Range("L2").Value = "stimuli/img_5jp4f.png"
$ws.Range("M2").Value = 84.85714285714286
$ws.Range("N2").Value = 67.83333333333333
$ws.Range("O2").Value = 76.3452380952381
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 9

# Row 3
$ws.Range("F3").Value = 363
$ws.Range("L3").Value = "stimuli/img_c0ecw.png"
$ws.Range("M3").Value = 18.88888888888889
$ws.Range("N3").Value = 17.82222222222222
$ws.Range("O3").Value = 18.35555555555555
$ws.Range("P3").Value = 45
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1

# Row 4
$ws.Range("F4").Value = 364
$ws.Range("J4").Value = "new"
$ws.Range("K4").Value = "f"
$ws.Range("L4").Value = "stimuli/img_koooi.png"
$ws.Range("M4").Value = 63.95454545454545
$ws.Range("N4").Value = 44.56818181818182
$ws.Range("O4").Value = 54.26136363636364
$ws.Range("P4").Value = 44
$ws.Range("Q4").Value = 5
$ws.Range("R4").Value = 5
$ws.Range("S4").Value = 5
$ws.Range("I4").ClearContents()

# Row 5
$ws.Range("F5").Value = 365
$ws.Range("J5").Value = "new"
$ws.Range("K5").Value = "f"
$ws.Range("L5").Value = "stimuli/img_jpjeg.png"
$ws.Range("M5").Value = 90.90697674418605
$ws.Range("N5").Value = 74.3953488372093
$ws.Range("O5").Value = 82.65116279069767
$ws.Range("P5").Value = 43
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 10
$ws.Range("I5").ClearContents()

# Row 6
$ws.Range("F6").Value = 366
$ws.Range("J6").Value = "new"
$ws.Range("K6").Value = "f"
$ws.Range("L6").Value = "stimuli/img_lpr0l.png"
$ws.Range("M6").Value = 77.04651162790698
$ws.Range("N6").Value = 59.86046511627907
$ws.Range("O6").Value = 68.45348837209303
$ws.Range("P6").Value = 43
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 7
$ws.Range("I6").ClearContents()

# Row 7
$ws.Range("F7").Value = 367
$ws.Range("L7").Value = "stimuli/img_syam3.png"
$ws.Range("M7").Value = 41.32432432432432
$ws.Range("N7").Value = 26.2972972972973
$ws.Range("O7").Value = 33.81081081081081
$ws.Range("P7").Value = 37
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 2

# Row 8
$ws.Range("F8").Value = 368
$ws.Range("J8").Value = "new"
$ws.Range("K8").Value = "f"
$ws.Range("L8").Value = "stimuli/img_yosqb.png"
$ws.Range("M8").Value = 50.88372093023256
$ws.Range("N8").Value = 30.11627906976744
$ws.Range("O8").Value = 40.5
$ws.Range("P8").Value = 43
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 3
$ws.Range("S8").Value = 3
$ws.Range("I8").ClearContents()

# Row 9
$ws.Range("F9").Value = 369
$ws.Range("I9").Value = "target"
$ws.Range("J9").Value = "old"
$ws.Range("K9").Value = "j"
$ws.Range("L9").Value = "stimuli/img_astid.png"
$ws.Range("M9").Value = 31.04761904761905
$ws.Range("N9").Value = 25.54761904761905
$ws.Range("O9").Value = 28.29761904761905
$ws.Range("P9").Value = 42
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2

# Row 10
$ws.Range("F10").Value = 370
$ws.Range("J10").Value = "new"
$ws.Range("K10").Value = "f"
$ws.Range("L10").Value = "stimuli/img_1iam9.png"
$ws.Range("M10").Value = 31.94594594594595
$ws.Range("N10").Value = 23.97297297297297
$ws.Range("O10").Value = 27.95945945945946
$ws.Range("P10").Value = 37
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 2
$ws.Range("S10").Value = 2
$ws.Range("I10").ClearContents()

# Row 11
$ws.Range("F11").Value = 371
$ws.Range("L11").Value = "stimuli/img_vnxft.png"
$ws.Range("M11").Value = 53.22727272727273
$ws.Range("N11").Value = 34.84090909090909
$ws.Range("O11").Value = 44.03409090909091
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = 3
$ws.Range("S11").Value = 3

# Row 12
$ws.Range("F12").Value = 372
$ws.Range("L12").Value = "stimuli/img_tujn3.png"
$ws.Range("M12").Value = 81.4090909090909
$ws.Range("N12").Value = 62.52272727272727
$ws.Range("O12").Value = 71.9659090909091
$ws.Range("P12").Value = 44

# Row 13
$ws.Range("F13").Value = 373
$ws.Range("L13").Value = "stimuli/img_vgh2g.png"
$ws.Range("M13").Value = 93.81395348837209
$ws.Range("N13").Value = 78.27906976744185
$ws.Range("O13").Value = 86.04651162790697
$ws.Range("P13").Value = 43
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 10

# Row 14
$ws.Range("F14").Value = 374
$ws.Range("J14").Value = "new"
$ws.Range("K14").Value = "f"
$ws.Range("L14").Value = "stimuli/img_tv8e2.png"
$ws.Range("M14").Value = 71.93023255813954
$ws.Range("N14").Value = 50.25581395348837
$ws.Range("O14").Value = 61.09302325581395
$ws.Range("P14").Value = 43
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 6
$ws.Range("S14").Value = 6
$ws.Range("I14").ClearContents()

# Row 15
$ws.Range("F15").Value = 375
$ws.Range("L15").Value = "stimuli/img_rg4in.png"
$ws.Range("M15").Value = 49.3695652173913
$ws.Range("N15").Value = 30.21739130434782
$ws.Range("O15").Value = 39.79347826086956
$ws.Range("P15").Value = 46
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 3

# Row 16
$ws.Range("F16").Value = 376
$ws.Range("I16").Value = "target"
$ws.Range("J16").Value = "old"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_1zhz6.png"
$ws.Range("M16").Value = 49.02272727272727
$ws.Range("N16").Value = 32.77272727272727
$ws.Range("O16").Value = 40.89772727272727
$ws.Range("P16").Value = 44
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 3

# Row 17
$ws.Range("F17").Value = 377
$ws.Range("L17").Value = "stimuli/img_of8d6.png"
$ws.Range("M17").Value = 26.04878048780488
$ws.Range("N17").Value = 19.14634146341463
$ws.Range("O17").Value = 22.59756097560975
$ws.Range("P17").Value = 41
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 1
$ws.Range("S17").Value = 1

# Row 18
$ws.Range("F18").Value = 378
$ws.Range("L18").Value = "stimuli/img_njhlh.png"
$ws.Range("M18").Value = 59.74418604651163
$ws.Range("N18").Value = 41.51162790697674
$ws.Range("O18").Value = 50.62790697674419

# Row 19
$ws.Range("F19").Value = 379
$ws.Range("H19").Value = "living_rooms"
$ws.Range("J19").Value = "new"
$ws.Range("L19").Value = "stimuli/img_g9od8.png"
$ws.Range("M19").Value = 59.34883720930232
$ws.Range("N19").Value = 37.83720930232558
$ws.Range("O19").Value = 48.59302325581395
$ws.Range("P19").Value = 43
$ws.Range("Q19").Value = 4
$ws.Range("R19").Value = 4
$ws.Range("S19").Value = 4

# Row 20
$ws.Range("F20").Value = 380
$ws.Range("I20").Value = "target"
$ws.Range("J20").Value = "old"
$ws.Range("K20").Value = "j"
$ws.Range("L20").Value = "stimuli/img_73pyk.png"
$ws.Range("M20").Value = 69.27659574468085
$ws.Range("N20").Value = 47.27659574468085
$ws.Range("O20").Value = 58.27659574468085
$ws.Range("Q20").Value = 5
$ws.Range("R20").Value = 5
$ws.Range("S20").Value = 5

# Row 21
$ws.Range("F21").Value = 381
$ws.Range("L21").Value = "stimuli/img_vh7v8.png"
$ws.Range("M21").Value = 78.70454545454545
$ws.Range("N21").Value = 59.63636363636363
$ws.Range("O21").Value = 69.17045454545455
$ws.Range("P21").Value = 44
$ws.Range("Q21").Value = 7
$ws.Range("R21").Value = 7
$ws.Range("S21").Value = 7

# Row 22
$ws.Range("F22").Value = 382
$ws.Range("J22").Value = "new"
$ws.Range("K22").Value = "f"
$ws.Range("L22").Value = "stimuli/img_c89x3.png"
$ws.Range("M22").Value = 72.8695652173913
$ws.Range("N22").Value = 49.65217391304348
$ws.Range("O22").Value = 61.26086956521739
$ws.Range("P22").Value = 46
$ws.Range("Q22").Value = 6
$ws.Range("R22").Value = 6
$ws.Range("S22").Value = 6
$ws.Range("I22").ClearContents()

# Row 23
$ws.Range("F23").Value = 383
$ws.Range("J23").Value = "catch"
$ws.Range("K23").Value = "f"
$ws.Range("L23").Value = "stimuli/catch_21.jpg"
$ws.Range("H23").ClearContents()
$ws.Range("I23").ClearContents()
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("O23").ClearContents()
$ws.Range("P23").ClearContents()
$ws.Range("Q23").ClearContents()
$ws.Range("R23").ClearContents()
$ws.Range("S23").ClearContents()

# Row 24
$ws.Range("F24").Value = 384
$ws.Range("L24").Value = "stimuli/img_3m61b.png"
$ws.Range("M24").Value = 81.97619047619048
$ws.Range("N24").Value = 63.23809523809524
$ws.Range("O24").Value = 72.60714285714286
$ws.Range("P24").Value = 42
$ws.Range("Q24").Value = 8
$ws.Range("R24").Value = 8
$ws.Range("S24").Value = 8

# Row 25
$ws.Range("F25").Value = 385
$ws.Range("J25").Value = "new"
$ws.Range("K25").Value = "f"
$ws.Range("L25").Value = "stimuli/img_pna7l.png"
$ws.Range("M25").Value = 85.53333333333333
$ws.Range("N25").Value = 67.97777777777777
$ws.Range("O25").Value = 76.75555555555556
$ws.Range("P25").Value = 45
$ws.Range("Q25").Value = 9
$ws.Range("R25").Value = 9
$ws.Range("S25").Value = 9
$ws.Range("I25").ClearContents()

# Row 26
$ws.Range("F26").Value = 386
$ws.Range("L26").Value = "stimuli/img_i6wsx.png"
$ws.Range("M26").Value = 79.07142857142857
$ws.Range("N26").Value = 58
$ws.Range("O26").Value = 68.53571428571428
$ws.Range("P26").Value = 42
$ws.Range("Q26").Value = 7
$ws.Range("R26").Value = 7
$ws.Range("S26").Value = 7

# Row 27
$ws.Range("F27").Value = 387
$ws.Range("J27").Value = "new"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_v771n.png"
$ws.Range("M27").Value = 22.325
$ws.Range("N27").Value = 15.25
$ws.Range("O27").Value = 18.7875
$ws.Range("P27").Value = 40
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
$ws.Range("I27").ClearContents()

# Row 28
$ws.Range("F28").Value = 388
$ws.Range("I28").Value = "target"
$ws.Range("J28").Value = "old"
$ws.Range("K28").Value = "j"
$ws.Range("L28").Value = "stimuli/img_tn8ys.png"
$ws.Range("M28").Value = 86.70454545454545
$ws.Range("N28").Value = 72.4090909090909
$ws.Range("O28").Value = 79.55681818181819
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = 10
$ws.Range("S28").Value = 10

# Row 29
$ws.Range("F29").Value = 389
$ws.Range("L29").Value = "stimuli/img_g13d5.png"
$ws.Range("M29").Value = 73
$ws.Range("N29").Value = 51.51111111111111
$ws.Range("O29").Value = 62.25555555555556
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6

# Row 30 (new trial row)
$ws.Range("A30").Value = 14
$ws.Range("B30").Value = "memory"
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 29
$ws.Range("F30").Value = 390
$ws.Range("G30").Value = "bedrooms"
$ws.Range("H30").Value = "living_rooms"
$ws.Range("I30").Value = "target"
$ws.Range("J30").Value = "old"
$ws.Range("K30").Value = "j"
$ws.Range("L30").Value = "stimuli/img_9oofc.png"
$ws.Range("M30").Value = 82.47619047619048
$ws.Range("N30").Value = 65.5
$ws.Range("O30").Value = 73.98809523809524
$ws.Range("P30").Value = 42
$ws.Range("Q30").Value = 8
$ws.Range("R30").Value = 8
$ws.Range("S30").Value = 8

# Row 31 (new trial row)
$ws.Range("A31").Value = 14
$ws.Range("B31").Value = "memory"
$ws.Range("C31").Value = 6
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 30
$ws.Range("F31").Value = 391
$ws.Range("G31").Value = "bedrooms"
$ws.Range("H31").Value = "living_rooms"
$ws.Range("J31").Value = "new"
$ws.Range("K31").Value = "f"
$ws.Range("L31").Value = "stimuli/img_x9w7o.png"
$ws.Range("M31").Value = 92.38888888888889
$ws.Range("N31").Value = 72.94444444444444
$ws.Range("O31").Value = 82.66666666666666
$ws.Range("P31").Value = 36
$ws.Range("Q31").Value = 10
$ws.Range("R31").Value = 10
$ws.Range("S31").Value = 10

# Row 32 (new trial row)
$ws.Range("A32").Value = 14
$ws.Range("B32").Value = "memory"
$ws.Range("C32").Value = 6
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 31
$ws.Range("F32").Value = 392
$ws.Range("G32").Value = "bedrooms"
$ws.Range("H32").Value = "living_rooms"
$ws.Range("I32").Value = "target"
$ws.Range("J32").Value = "old"
$ws.Range("K32").Value = "j"
$ws.Range("L32").Value = "stimuli/img_swq34.png"
$ws.Range("M32").Value = 64.11363636363636
$ws.Range("N32").Value = 43.04545454545455
$ws.Range("O32").Value = 53.57954545454545
$ws.Range("P32").Value = 44
$ws.Range("Q32").Value = 5
$ws.Range("R32").Value = 5
$ws.Range("S32").Value = 5

# Row 33 (new trial row)
$ws.Range("A33").Value = 14
$ws.Range("B33").Value = "memory"
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 32
$ws.Range("F33").Value = 393
$ws.Range("G33").Value = "bedrooms"
$ws.Range("H33").Value = "living_rooms"
$ws.Range("J33").Value = "new"
$ws.Range("K33").Value = "f"
$ws.Range("L33").Value = "stimuli/img_165pk.png"
$ws.Range("M33").Value = 85.73333333333333
$ws.Range("N33").Value = 69.22222222222223
$ws.Range("O33").Value = 77.47777777777779
$ws.Range("P33").Value = 45
$ws.Range("Q33").Value = 9
$ws.Range("R33").Value = 9
$ws.Range("S33").Value = 9

# Row 34 (new trial row)
$ws.Range("A34").Value = 14
$ws.Range("B34").Value = "memory"
$ws.Range("C34").Value = 6
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 33
$ws.Range("F34").Value = 394
$ws.Range("G34").Value = "bedrooms"
$ws.Range("H34").Value = "living_rooms"
$ws.Range("I34").Value = "target"
$ws.Range("J34").Value = "old"
$ws.Range("K34").Value = "j"
$ws.Range("L34").Value = "stimuli/img_xr3up.png"
$ws.Range("M34").Value = 76.24444444444444
$ws.Range("N34").Value = 55.88888888888889
$ws.Range("O34").Value = 66.06666666666666
$ws.Range("P34").Value = 45
$ws.Range("Q34").Value = 7
$ws.Range("R34").Value = 7
$ws.Range("S34").Value = 7

# Row 35 (new trial row)
$ws.Range("A35").Value = 14
$ws.Range("B35").Value = "memory"
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 34
$ws.Range("F35").Value = 395
$ws.Range("G35").Value = "bedrooms"
$ws.Range("H35").Value = "living_rooms"
$ws.Range("I35").Value = "target"
$ws.Range("J35").Value = "old"
$ws.Range("K35").Value = "j"
$ws.Range("L35").Value = "stimuli/img_j4ttn.png"
$ws.Range("M35").Value = 12.61904761904762
$ws.Range("N35").Value = 11.42857142857143
$ws.Range("O35").Value = 12.02380952380952
$ws.Range("P35").Value = 42
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = 1
$ws.Range("S35").Value = 1

# Row 36 (new trial row)
$ws.Range("A36").Value = 14
$ws.Range("B36").Value = "memory"
$ws.Range("C36").Value = 6
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 35
$ws.Range("F36").Value = 396
$ws.Range("G36").Value = "bedrooms"
$ws.Range("H36").Value = "living_rooms"
$ws.Range("I36").Value = "target"
$ws.Range("J36").Value = "old"
$ws.Range("K36").Value = "j"
$ws.Range("L36").Value = "stimuli/img_rru0v.png"
$ws.Range("M36").Value = 56.45238095238095
$ws.Range("N36").Value = 39.42857142857143
$ws.Range("O36").Value = 47.94047619047619
$ws.Range("P36").Value = 42
$ws.Range("Q36").Value = 4
$ws.Range("R36").Value = 4
$ws.Range("S36").Value = 4

# Row 37 (new trial row)
$ws.Range("A37").Value = 14
$ws.Range("B37").Value = "memory"
$ws.Range("C37").Value = 6
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = 36
$ws.Range("F37").Value = 397
$ws.Range("G37").Value = "bedrooms"
$ws.Range("H37").Value = "living_rooms"
$ws.Range("I37").Value = "target"
$ws.Range("J37").Value = "old"
$ws.Range("K37").Value = "j"
$ws.Range("L37").Value = "stimuli/img_ra2nm.png"
$ws.Range("M37").Value = 70.75
$ws.Range("N37").Value = 50.375
$ws.Range("O37").Value = 60.5625
$ws.Range("P37").Value = 40
$ws.Range("Q37").Value = 6
$ws.Range("R37").Value = 6
$ws.Range("S37").Value = 6

# Row 38 (new trial row)
$ws.Range("A38").Value = 14
$ws.Range("B38").Value = "memory"
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 2
$ws.Range("E38").Value = 37
$ws.Range("F38").Value = 398
$ws.Range("G38").Value = "bedrooms"
$ws.Range("H38").Value = "living_rooms"
$ws.Range("I38").Value = "target"
$ws.Range("J38").Value = "old"
$ws.Range("K38").Value = "j"
$ws.Range("L38").Value = "stimuli/img_ac0ey.png"
$ws.Range("M38").Value = 86.62222222222222
$ws.Range("N38").Value = 70.02222222222223
$ws.Range("O38").Value = 78.32222222222222
$ws.Range("P38").Value = 45
$ws.Range("Q38").Value = 9
$ws.Range("R38").Value = 9
$ws.Range("S38").Value = 9

# Row 39 (new trial row)
$ws.Range("A39").Value = 14
$ws.Range("B39").Value = "memory"
$ws.Range("C39").Value = 6
$ws.Range("D39").Value = 2
$ws.Range("E39").Value = 38
$ws.Range("F39").Value = 399
$ws.Range("G39").Value = "bedrooms"
$ws.Range("H39").Value = "living_rooms"
$ws.Range("I39").Value = "target"
$ws.Range("J39").Value = "old"
$ws.Range("K39").Value = "j"
$ws.Range("L39").Value = "stimuli/img_rych7.png"
$ws.Range("M39").Value = 30.4468085106383
$ws.Range("N39").Value = 23.4468085106383
$ws.Range("O39").Value = 26.9468085106383
$ws.Range("P39").Value = 47
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = 2
$ws.Range("S39").Value = 2

# Row 40 (new trial row)
$ws.Range("A40").Value = 14
$ws.Range("B40").Value = "memory"
$ws.Range("C40").Value = 6
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = 39
$ws.Range("F40").Value = 400
$ws.Range("G40").Value = "bedrooms"
$ws.Range("H40").Value = "living_rooms"
$ws.Range("J40").Value = "new"
$ws.Range("K40").Value = "f"
$ws.Range("L40").Value = "stimuli/img_jkm86.png"
$ws.Range("M40").Value = 58.32558139534883
$ws.Range("N40").Value = 38.65116279069768
$ws.Range("O40").Value = 48.48837209302326
$ws.Range("P40").Value = 43
$ws.Range("Q40").Value = 4
$ws.Range("R40").Value = 4
$ws.Range("S40").Value = 4

# Row 41 (new trial row)
$ws.Range("A41").Value = 14
$ws.Range("B41").Value = "memory"
$ws.Range("C41").Value = 6
$ws.Range("D41").Value = 2
$ws.Range("E41").Value = 40
$ws.Range("F41").Value = 401
$ws.Range("G41").Value = "bedrooms"
$ws.Range("H41").Value = "living_rooms"
$ws.Range("J41").Value = "new"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_xpco9.png"
$ws.Range("M41").Value = 81.55555555555556
$ws.Range("N41").Value = 64.68888888888888
$ws.Range("O41").Value = 73.12222222222222
$ws.Range("P41").Value = 45
$ws.Range("Q41").Value = 8
$ws.Range("R41").Value = 8
$ws.Range("S41").Value = 8
